$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append (date serial, nuovi pos., somma mobile 7gg., per 100mila ab.)
$data = @(
    @(44432, 1, 20, 157.0475068708284),
    @(44433, 0, 19, 149.195131527287),
    @(44434, 0, 17, 133.4903808402042),
    @(44435, 4, 20, 157.0475068708284),
    @(44436, 3, 11, 86.37612877895563),
    @(44437, 1, 11, 86.37612877895563),
    @(44438, 2, 11, 86.37612877895563),
    @(44439, 0, 10, 78.52375343541422),
    @(44440, 0, 10, 78.52375343541422)
)

$startRow = 358

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    # Match formatting of the preceding data row first (column A carries
    # the date style/border, B-D stay default like the rest of the table).
    $ws.Cells.Item($r - 1, 1).Copy() | Out-Null
    $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
$excel.CutCopyMode = 0
